# Apply the "Promotion/Demotion" menu entry to the HR sheet of the MenuList
# workbook, mirroring the target OOXML diff:
#  - sheet2.xml ("HR"): append a new row (row 66) with the Promotion/Demotion
#    menu item, extend dimension/selection accordingly.
#  - sharedStrings.xml: the two new literal strings get appended implicitly
#    as a side effect of writing the new cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HR")

# New row to append right after the existing last row (65).
$newRow = 66

# Write the URL (column D) before the menu name (column A) so that the
# shared-string table receives "/hr/promotion-demotion/" before
# "Promotion/Demotion", matching the order new unique strings were added
# in the source workbook.
$ws.Range("D" + $newRow).Value = "/hr/promotion-demotion/"
$ws.Range("A" + $newRow).Value = "Promotion/Demotion"
$ws.Range("B" + $newRow).Value = "Yes"
$ws.Range("C" + $newRow).Value = "Management"
$ws.Range("E" + $newRow).Value = "fas fa-users"
$ws.Range("F" + $newRow).Value = 26

# Update the sheet view: select the cell below the freshly added row and
# scroll so that the new row is visible, matching the author's view state
# after adding the entry.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A67").Select() | Out-Null
